$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Force the cell to stay a text value (avoid Excel auto-converting
    # numeric-looking strings like "1.00" or "551.15" into numbers, which
    # would drop trailing zeros / change the stored type), then restore the
    # default "Normal" style so no stray formatting/style index is left on
    # the cell.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Cells.Item(2, 4) "67.770.48"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.96%  "

# Row 3 - Ethereum
Set-TextValue $ws.Cells.Item(3, 4) "2.425.27"
Set-TextValue $ws.Cells.Item(3, 5) "  -2.27%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Cells.Item(4, 5) "  +0.07%  "

# Row 5 - BNB
Set-TextValue $ws.Cells.Item(5, 4) "551.15"
Set-TextValue $ws.Cells.Item(5, 5) "  -2.37%  "

# Row 6 - Solana
Set-TextValue $ws.Cells.Item(6, 4) "159.45"
Set-TextValue $ws.Cells.Item(6, 5) "  -2.34%  "

# Row 7 - USDC
Set-TextValue $ws.Cells.Item(7, 5) "  -0.01%  "

# Row 8 - XRP
Set-TextValue $ws.Cells.Item(8, 5) "  -2.65%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Cells.Item(9, 4) "2.428.13"
Set-TextValue $ws.Cells.Item(9, 5) "  -2.05%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Cells.Item(10, 4) "0.145"
Set-TextValue $ws.Cells.Item(10, 5) "  -7.50%  "

# Row 11 - TRON
Set-TextValue $ws.Cells.Item(11, 5) "  -1.87%  "

# Row 12 - Cardano
Set-TextValue $ws.Cells.Item(12, 5) "  -5.68%  "

# Row 13 - Toncoin
Set-TextValue $ws.Cells.Item(13, 4) "4.69"
Set-TextValue $ws.Cells.Item(13, 5) "  -3.96%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Cells.Item(14, 4) "2.875.38"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.11%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Cells.Item(15, 4) "67.646.66"
Set-TextValue $ws.Cells.Item(15, 5) "  -2.29%  "

# Row 16 - ShibaInu
Set-TextValue $ws.Cells.Item(16, 5) "  -6.09%  "

# Row 17 - Avalanche
Set-TextValue $ws.Cells.Item(17, 4) "22.74"
Set-TextValue $ws.Cells.Item(17, 5) "  -5.91%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Cells.Item(18, 4) "2.412.89"
Set-TextValue $ws.Cells.Item(18, 5) "  -2.63%  "

# Row 19 - Chainlink
Set-TextValue $ws.Cells.Item(19, 4) "10.65"
Set-TextValue $ws.Cells.Item(19, 5) "  -4.32%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Cells.Item(20, 4) "336.13"
Set-TextValue $ws.Cells.Item(20, 5) "  -2.28%  "

# Row 21 - Uniswap
Set-TextValue $ws.Cells.Item(21, 5) "  -5.30%  "

# Row 22 - Polkadot
Set-TextValue $ws.Cells.Item(22, 5) "  -4.10%  "

# Row 23 - Dai
Set-TextValue $ws.Cells.Item(23, 4) "0.998"
Set-TextValue $ws.Cells.Item(23, 5) "  -0.15%  "

# Row 24 - SuiNetwork
Set-TextValue $ws.Cells.Item(24, 4) "1.80"
Set-TextValue $ws.Cells.Item(24, 5) "  -5.47%  "

# Row 25 - Litecoin
Set-TextValue $ws.Cells.Item(25, 4) "65.91"
Set-TextValue $ws.Cells.Item(25, 5) "  -4.71%  "

# Row 26 - WrappedeETH
Set-TextValue $ws.Cells.Item(26, 4) "2.554.71"
Set-TextValue $ws.Cells.Item(26, 5) "  -2.03%  "

# Row 27 - NEARProtocol
Set-TextValue $ws.Cells.Item(27, 5) "  -7.34%  "

# Row 28 - Binance-PegBSC-USD
Set-TextValue $ws.Cells.Item(28, 5) "  +0.75%  "

# Row 29 - Aptos
Set-TextValue $ws.Cells.Item(29, 4) "7.92"
Set-TextValue $ws.Cells.Item(29, 5) "  -8.32%  "

# Row 30 - PEPE
Set-TextValue $ws.Cells.Item(30, 4) "0.0₃0804"
Set-TextValue $ws.Cells.Item(30, 5) "  -7.38%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue $ws.Cells.Item(31, 4) "6.97"
Set-TextValue $ws.Cells.Item(31, 5) "  -8.88%  "

# Row 32 - FirstDigitalUSD
Set-TextValue $ws.Cells.Item(32, 4) "1.00"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.23%  "

# Row 33 - Bittensor
Set-TextValue $ws.Cells.Item(33, 4) "409.83"
Set-TextValue $ws.Cells.Item(33, 5) "  -6.89%  "

# Row 34 - PancakeSwap
Set-TextValue $ws.Cells.Item(34, 4) "1.60"
Set-TextValue $ws.Cells.Item(34, 5) "  -5.38%  "

# Row 35 - Fetch.AI
Set-TextValue $ws.Cells.Item(35, 5) "  -6.64%  "

# Row 36 - Monero
Set-TextValue $ws.Cells.Item(36, 4) "157.30"
Set-TextValue $ws.Cells.Item(36, 5) "  +1.30%  "

# Row 37 - WhiteBITCoin
Set-TextValue $ws.Cells.Item(37, 4) "18.97"
Set-TextValue $ws.Cells.Item(37, 5) "  -0.18%  "

# Row 39 - Kaspa
Set-TextValue $ws.Cells.Item(39, 5) "  -5.51%  "

# Row 40 - EthereumClassic
Set-TextValue $ws.Cells.Item(40, 4) "17.55"
Set-TextValue $ws.Cells.Item(40, 5) "  -2.86%  "

# Row 41 - PolygonEcosystemToken
Set-TextValue $ws.Cells.Item(41, 4) "0.296"
Set-TextValue $ws.Cells.Item(41, 5) "  -5.12%  "

# Row 42 - RenderToken
Set-TextValue $ws.Cells.Item(42, 5) "  -7.08%  "

# Row 43 - Stacks
Set-TextValue $ws.Cells.Item(43, 4) "1.44"
Set-TextValue $ws.Cells.Item(43, 5) "  -7.89%  "

# Row 44 - ImmutableX
Set-TextValue $ws.Cells.Item(44, 4) "1.05"
Set-TextValue $ws.Cells.Item(44, 5) "  -1.37%  "

# Row 45 - Aave
Set-TextValue $ws.Cells.Item(45, 4) "132.02"

# Row 46 - dogwifhat
Set-TextValue $ws.Cells.Item(46, 4) "1.99"
Set-TextValue $ws.Cells.Item(46, 5) "  -7.77%  "

# Row 47 - Filecoin
Set-TextValue $ws.Cells.Item(47, 4) "3.26"
Set-TextValue $ws.Cells.Item(47, 5) "  -4.57%  "

# Row 48 - Cronos
Set-TextValue $ws.Cells.Item(48, 4) "0.0708"
Set-TextValue $ws.Cells.Item(48, 5) "  -2.02%  "

# Row 49 - was ARBITRUM, now Mantle (rows 49/50 swapped content)
Set-TextValue $ws.Cells.Item(49, 2) "Mantle"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Cells.Item(49, 4) "0.550"
Set-TextValue $ws.Cells.Item(49, 5) "  -3.33%  "

# Row 50 - was Mantle, now ARBITRUM
Set-TextValue $ws.Cells.Item(50, 2) "ARBITRUM"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(50, 4) "0.465"
Set-TextValue $ws.Cells.Item(50, 5) "  -8.67%  "

# Row 51 - Stellar
Set-TextValue $ws.Cells.Item(51, 5) "  -2.13%  "
